$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 142 (pushes the old row 142 and below down by one).
# This automatically shifts formulas (C144, E144/E145, G144/G145, C146, etc.)
# and extends dataValidation/merged ranges that span the insertion point.
$ws.Rows.Item(142).Insert()

# The freshly inserted row 142 copied formatting straight down from row 141,
# but two columns (D and J) actually ended up with a different style in the
# real edit (picked up from elsewhere), so fix those two cells' formats by
# copying from cells that already carry the desired style.
$ws.Range("D114").Copy() | Out-Null
$ws.Range("D142").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("K141").Copy() | Out-Null
$ws.Range("J142").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new task row content.
$ws.Cells.Item(142, 1).Value2 = 22
$ws.Cells.Item(142, 2).Value2 = "Interface Design"
$ws.Cells.Item(142, 3).Value2 = "MockUps"
$ws.Cells.Item(142, 4).Value2 = "[TASK]"
$ws.Cells.Item(142, 5).Value2 = "Umbau von Ordner Struktur und anlegen von Tablets"
$ws.Cells.Item(142, 6).Value2 = 44464
$ws.Cells.Item(142, 7).Value2 = 44481
$ws.Range("I142").Formula = "=ROUNDUP(((SUM(K142-J142)*24*60/60)/0.25),0)*0.25"
$ws.Cells.Item(142, 10).Value2 = 0.55208333333333337
$ws.Cells.Item(142, 11).Value2 = 0.61458333333333337

# Row 142 now holds real data instead of being an empty placeholder row, so
# it should no longer be part of the Prefix dropdown validation range (the
# row insert auto-extended that validation's sqref to include it).
$ws.Range("D142").Validation.Delete()
